$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Product rows (written first so shared strings line up with target order) ---
$ws.Range("C2").Value = "banh  bao;12000;23"
$ws.Range("D2").Value = "banh  bao 2;12000;23"
$ws.Range("E2").Value = "banh  bao 3;23000;10"
$ws.Range("C3").Value = "thung rac;12000;23"
$ws.Range("D3").Value = "thung rac;50000;12"

# --- Customer cells ---
$ws.Range("A2").Value = "Customer 2"
$ws.Range("A3").Value = "Customer 2"

# --- Note cell (C1) ---
$ws.Range("C1").Value = "Note:  From C column to the right (end: V columm), each cell is the product data in the order via format <name>;<price>;<quantity> "
$ws.Range("C1").Font.Color = 255

# --- Customer cell font (distinguishable "family 2" font, closest achievable match) ---
$ws.Range("A2").Font.ThemeFont = 1
$ws.Range("A3").Font.ThemeFont = 1
$ws.Range("D3").Font.ThemeFont = 1

# --- Order date cells ---
$ws.Range("B2").Value = 44879
$ws.Range("B3").Value = 44879
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Range("B3").NumberFormat = "mm-dd-yy"

$ws.Range("C1").Characters(8,5).Font.ColorIndex = -4105
$ws.Range("C1").Characters(13,1).Font.Bold = $true
$ws.Range("C1").Characters(13,1).Font.ColorIndex = -4105
$ws.Range("C1").Characters(14,27).Font.ColorIndex = -4105
$ws.Range("C1").Characters(41,1).Font.Bold = $true
$ws.Range("C1").Characters(41,1).Font.ColorIndex = -4105
$ws.Range("C1").Characters(42,9).Font.ColorIndex = -4105

$ws.Range("C1").Characters(51,55).Font.Color = 0
$ws.Range("C1").Characters(106,6).Font.Color = 15773696
$ws.Range("C1").Characters(112,1).Font.Bold = $true
$ws.Range("C1").Characters(112,1).Font.Color = 0
$ws.Range("C1").Characters(113,7).Font.Color = 15773696
$ws.Range("C1").Characters(120,1).Font.Bold = $true
$ws.Range("C1").Characters(120,1).Font.Color = 0
$ws.Range("C1").Characters(121,11).Font.Color = 15773696

# --- Column width for column E ---
$ws.Columns.Item(5).ColumnWidth = 18.5

# --- View state ---
$ws.Range("C3").Select()

# --- Page setup ---
$ws.PageSetup.Orientation = 1
